## Optimising for healthy children, new version of ASD
# Applies the cell/formula/comment/view changes described by the target diff
# to the "Optimal funding scenario" sheet (and related workbook-level view
# state) of the Tanzania optimisationBudgets workbook.

$wb = $excel.ActiveWorkbook

$wsCurrent = $wb.Worksheets.Item("Current expenditure")
$wsOpt     = $wb.Worksheets.Item("Optimal funding scenario")

# --- Scenario 1 row (row 2) now mirrors the "free, programatically
#     optimised" pattern already used by rows 4/5: update the scenario
#     description, give it the averaged-spend formula (same formula/style
#     as D4/D5), and move the "include in analysis" check mark from E4
#     up to E2. ---

# C2: "N/A" -> "Free, programatically optimised" (same text already used in C4/C5)
$wsOpt.Range("C2").Value = $wsOpt.Range("C4").Value2

# D2: copy D4's number format / fill (style 13) onto D2, then give it the
# same averaging formula so the computed value matches (32,833,333.33...)
$wsOpt.Range("D4").Copy()
$wsOpt.Range("D2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$wsOpt.Range("D2").Formula = "=(20+30+45+45+35+22)/6 * 1000000"

# E2 gains the "x" check mark, E4 loses it
$wsOpt.Range("E2").Value = "x"
$wsOpt.Range("E4").ClearContents()

# --- New comment on D2, matching the wording already used on D4/D5 ---
$wsOpt.Range("D2").AddComment("Sam:" + [char]10 + "Average annual spending 2019-2024 from Yi-Kyoung")

# --- View state: "Optimal funding scenario" becomes the active / selected
#     tab (was "Current expenditure"), with the cursor on E2. ---
$wsCurrent.Range("B14").Select()
$wsOpt.Activate()
$wsOpt.Range("E2").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
